$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = '[''episode'', ''mulder'', ''doctor'', ''scully'', ''viewers'', ''episodes'', ''series'', ''television'', ''character'', ''season'', ''watched'', ''dwight'', ''broadcast'', ''trek'', ''jack'']'
$ws.Range("C2").Value2 = 0
$ws.Range("D2").Value2 = 0
$ws.Range("E2").Value2 = 0
$ws.Range("F2").Value2 = 0
$ws.Range("G2").Value2 = 0

$ws.Range("B3").Value2 = '[''album'', ''song'', ''madonna'', ''chart'', ''video'', ''music'', ''harrison'', ''songs'', ''carey'', ''band'', ''pop'', ''track'', ''recording'', ''billboard'', ''number'']'
$ws.Range("C3").Value2 = 0.7797157170189593
$ws.Range("D3").Value2 = 0.6961537650392843
$ws.Range("E3").Value2 = 0.6599544747492919
$ws.Range("F3").Value2 = 0.6585169491760476
$ws.Range("G3").Value2 = 0.6467437531982382

$ws.Range("B4").Value2 = '[''ship'', ''guns'', ''ships'', ''tons'', ''torpedo'', ''knots'', ''inch'', ''cruiser'', ''fleet'', ''gun'', ''deck'', ''admiral'', ''german'', ''turrets'', ''cruisers'']'
$ws.Range("C4").Value2 = 0.8995174911302904
$ws.Range("D4").Value2 = 0.8973708912050163
$ws.Range("E4").Value2 = 0.8949955454517495
$ws.Range("F4").Value2 = 0.894203330176475
$ws.Range("G4").Value2 = 0.8929476064634039

$ws.Range("B5").Value2 = '[''highway'', ''route'', ''road'', ''freeway'', ''interchange'', ''intersection'', ''terminus'', ''north'', ''east'', ''lane'', ''continues'', ''state'', ''avenue'', ''passes'', ''traffic'']'
$ws.Range("C5").Value2 = 0.8994529065196937
$ws.Range("D5").Value2 = 0.8979494389191461
$ws.Range("E5").Value2 = 0.8973824479064191
$ws.Range("F5").Value2 = 0.8960262822473475
$ws.Range("G5").Value2 = 0.8913213516555959

$ws.Range("B6").Value2 = '[''election'', ''hitler'', ''campaign'', ''party'', ''bush'', ''republican'', ''political'', ''vote'', ''labour'', ''president'', ''presidential'', ''democratic'', ''war'', ''government'', ''senate'']'
$ws.Range("C6").Value2 = 0
$ws.Range("D6").Value2 = 0
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = 0
$ws.Range("G6").Value2 = 0

$ws.Range("B7").Value2 = '[''tropical'', ''storm'', ''hurricane'', ''winds'', ''depression'', ''cyclone'', ''mph'', ''rainfall'', ''damage'', ''landfall'', ''wind'', ''utc'', ''flooding'', ''weakened'', ''intensity'']'
$ws.Range("C7").Value2 = 0.8994931954472305
$ws.Range("D7").Value2 = 0.8989416131641279
$ws.Range("E7").Value2 = 0.8983367875255148
$ws.Range("F7").Value2 = 0.8911209008282481
$ws.Range("G7").Value2 = 0.8898994247877255

$ws.Range("B8").Value2 = '[''king'', ''government'', ''polish'', ''political'', ''emperor'', ''son'', ''arab'', ''death'', ''military'', ''army'', ''died'', ''pope'', ''poland'', ''reign'', ''byzantine'']'
$ws.Range("C8").Value2 = 0
$ws.Range("D8").Value2 = 0
$ws.Range("E8").Value2 = 0
$ws.Range("F8").Value2 = 0
$ws.Range("G8").Value2 = 0

$ws.Range("B9").Value2 = '[''species'', ''shark'', ''genus'', ''sharks'', ''females'', ''prey'', ''eggs'', ''males'', ''cap'', ''birds'', ''nest'', ''habitat'', ''fin'', ''brown'', ''stem'']'
$ws.Range("C9").Value2 = 0.6708850893472472
$ws.Range("D9").Value2 = 0.6629392884877827
$ws.Range("E9").Value2 = 0.6068501438223768
$ws.Range("F9").Value2 = 0.5017763141975303
$ws.Range("G9").Value2 = 0

$ws.Range("B10").Value2 = '[''game'', ''player'', ''gameplay'', ''games'', ''players'', ''soundtrack'', ''mario'', ''playstation'', ''graphics'', ''mode'', ''hero'', ''nintendo'', ''characters'', ''released'', ''version'']'
$ws.Range("C10").Value2 = 0
$ws.Range("D10").Value2 = 0
$ws.Range("E10").Value2 = 0
$ws.Range("F10").Value2 = 0
$ws.Range("G10").Value2 = 0

$ws.Range("B11").Value2 = '[''innings'', ''runs'', ''league'', ''baseball'', ''nba'', ''batting'', ''season'', ''career'', ''team'', ''scored'', ''wickets'', ''basketball'', ''games'', ''rebounds'', ''batted'']'
$ws.Range("C11").Value2 = 0.8901149788288089
$ws.Range("D11").Value2 = 0.8710989796136231
$ws.Range("E11").Value2 = 0.8391870939273561
$ws.Range("F11").Value2 = 0.8278851306559769
$ws.Range("G11").Value2 = 0.8240034880070073

$ws.Range("B12").Value2 = '[''film'', ''films'', ''bond'', ''disney'', ''role'', ''movie'', ''actor'', ''starred'', ''cast'', ''box'', ''tamil'', ''director'', ''filming'', ''grossing'', ''production'']'
$ws.Range("C12").Value2 = 0
$ws.Range("D12").Value2 = 0
$ws.Range("E12").Value2 = 0
$ws.Range("F12").Value2 = 0
$ws.Range("G12").Value2 = 0

$ws.Range("B13").Value2 = '[''french'', ''british'', ''governor'', ''militia'', ''troops'', ''battle'', ''massachusetts'', ''boston'', ''command'', ''kentucky'', ''fort'', ''men'', ''army'', ''fleet'', ''general'']'
$ws.Range("C13").Value2 = 0
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 0

$ws.Range("B14").Value2 = '[''castle'', ''century'', ''station'', ''church'', ''castles'', ''built'', ''bridge'', ''tower'', ''pier'', ''railway'', ''bailey'', ''river'', ''building'', ''trains'', ''bristol'']'
$ws.Range("C14").Value2 = 0.8898604909005957
$ws.Range("D14").Value2 = 0.8868397759459766
$ws.Range("E14").Value2 = 0.8754026800263479
$ws.Range("F14").Value2 = 0.8653366724333088
$ws.Range("G14").Value2 = 0.8573847899841381

$ws.Range("B15").Value2 = '[''match'', ''wrestling'', ''championship'', ''wwe'', ''tag'', ''raw'', ''ring'', ''smackdown'', ''defeated'', ''heavyweight'', ''feud'', ''event'', ''wwf'', ''michaels'', ''title'']'
$ws.Range("C15").Value2 = 0.8994461037061843
$ws.Range("D15").Value2 = 0.8970629673216594
$ws.Range("E15").Value2 = 0.8929929559887967
$ws.Range("F15").Value2 = 0.8918549569086228
$ws.Range("G15").Value2 = 0.8901445872411098

$ws.Range("B16").Value2 = '[''building'', ''city'', ''park'', ''chicago'', ''memorial'', ''memorials'', ''library'', ''fountain'', ''indiana'', ''square'', ''galveston'', ''street'', ''buildings'', ''bay'', ''courthouse'']'
$ws.Range("C16").Value2 = 0
$ws.Range("D16").Value2 = 0
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 0
$ws.Range("G16").Value2 = 0

$ws.Range("B17").Value2 = '[''simpsons'', ''episodes'', ''episode'', ''homer'', ''stan'', ''parker'', ''kenny'', ''animated'', ''kyle'', ''nickelodeon'', ''voice'', ''lisa'', ''jake'', ''television'', ''voiced'']'
$ws.Range("C17").Value2 = 0.8410681958333622
$ws.Range("D17").Value2 = 0.8408707414726451
$ws.Range("E17").Value2 = 0.8342199186176124
$ws.Range("F17").Value2 = 0.789329704457635
$ws.Range("G17").Value2 = 0.7601764068907592

$ws.Range("B18").Value2 = '[''battalion'', ''brigade'', ''division'', ''regiment'', ''infantry'', ''battalions'', ''squadron'', ''training'', ''unit'', ''units'', ''raaf'', ''australian'', ''2nd'', ''japanese'', ''1st'']'
$ws.Range("C18").Value2 = 0.8976577662524322
$ws.Range("D18").Value2 = 0.8907184988740995
$ws.Range("E18").Value2 = 0.8847848941949625
$ws.Range("F18").Value2 = 0.8843924798218827
$ws.Range("G18").Value2 = 0.8779918499766377

$ws.Range("B19").Value2 = '[''olympics'', ''athletes'', ''olympic'', ''medal'', ''meter'', ''freestyle'', ''championships'', ''beijing'', ''medals'', ''seconds'', ''relay'', ''gold'', ''competed'', ''games'', ''summer'']'
$ws.Range("C19").Value2 = 0.8995757378356433
$ws.Range("D19").Value2 = 0.8992369611502073
$ws.Range("E19").Value2 = 0.8962133967040117
$ws.Range("F19").Value2 = 0.8959980814952551
$ws.Range("G19").Value2 = 0.8939842142569264

$ws.Range("B20").Value2 = '[''oxford'', ''cambridge'', ''race'', ''boat'', ''blues'', ''rowed'', ''rowing'', ''lengths'', ''crews'', ''rower'', ''thames'', ''races'', ''crew'', ''universities'', ''umpired'']'
$ws.Range("C20").Value2 = 0.8933404872102932
$ws.Range("D20").Value2 = 0.8931113782617631
$ws.Range("E20").Value2 = 0.885902464354956
$ws.Range("F20").Value2 = 0.8857321694324698
$ws.Range("G20").Value2 = 0.8841812799749215

$ws.Range("B21").Value2 = '[''poem'', ''poems'', ''poetry'', ''poet'', ''han'', ''shakespeare'', ''ode'', ''riley'', ''text'', ''sanskrit'', ''texts'', ''poetic'', ''literary'', ''works'', ''smart'']'
$ws.Range("C21").Value2 = 0.7934033862609332
$ws.Range("D21").Value2 = 0.7498369472030213
$ws.Range("E21").Value2 = 0.7379377015299917
$ws.Range("F21").Value2 = 0.6746928011070878
$ws.Range("G21").Value2 = 0.6746928011070878

$ws.Range("B22").Value2 = '[''breed'', ''horses'', ''horse'', ''breeds'', ''dog'', ''dogs'', ''stud'', ''breeding'', ''bred'', ''arabian'', ''stallion'', ''riding'', ''breeders'', ''pony'', ''stakes'']'
$ws.Range("C22").Value2 = 0.8816342584032073
$ws.Range("D22").Value2 = 0.8789928337046216
$ws.Range("E22").Value2 = 0.8704516537222983
$ws.Range("F22").Value2 = 0.8349591423823824
$ws.Range("G22").Value2 = 0.8335245520367489

$ws.Range("B23").Value2 = '[''creek'', ''watershed'', ''dam'', ''volcano'', ''lava'', ''pipeline'', ''flows'', ''volcanic'', ''trout'', ''mountain'', ''park'', ''eruption'', ''river'', ''feet'', ''cubic'']'
$ws.Range("C23").Value2 = 0.8962382467061845
$ws.Range("D23").Value2 = 0.8748386294659284
$ws.Range("E23").Value2 = 0.867797248506632
$ws.Range("F23").Value2 = 0.8564198104348976
$ws.Range("G23").Value2 = 0.8476569538688893

$ws.Range("B24").Value2 = '[''aircraft'', ''engine'', ''flight'', ''fuselage'', ''wing'', ''air'', ''engines'', ''prototype'', ''radar'', ''fighter'', ''fuel'', ''raf'', ''speed'', ''car'', ''testing'']'
$ws.Range("C24").Value2 = 0.8923838961112569
$ws.Range("D24").Value2 = 0.8550846709305449
$ws.Range("E24").Value2 = 0.8372260935178961
$ws.Range("F24").Value2 = 0.8259438922902133
$ws.Range("G24").Value2 = 0.8191061185877834

$ws.Range("B25").Value2 = '[''trains'', ''locomotives'', ''locomotive'', ''oslo'', ''line'', ''nok'', ''train'', ''station'', ''tunnel'', ''railway'', ''class'', ''rail'', ''railways'', ''passenger'', ''trondheim'']'
$ws.Range("C25").Value2 = 0.8999502951522912
$ws.Range("D25").Value2 = 0.8961412141544224
$ws.Range("E25").Value2 = 0.8946205474495311
$ws.Range("F25").Value2 = 0.894024604961298
$ws.Range("G25").Value2 = 0.8928936403254673

$ws.Range("B26").Value2 = '[''manga'', ''anime'', ''comics'', ''stories'', ''story'', ''comic'', ''magazine'', ''volume'', ''volumes'', ''fiction'', ''gay'', ''characters'', ''pulp'', ''published'', ''issue'']'
$ws.Range("C26").Value2 = 0.8971397471299966
$ws.Range("D26").Value2 = 0.8950281123517705
$ws.Range("E26").Value2 = 0.887937349280863
$ws.Range("F26").Value2 = 0.8853554957904606
$ws.Range("G26").Value2 = 0.8776313937324671

$ws.Range("B27").Value2 = '[''lap'', ''race'', ''drivers'', ''laps'', ''pit'', ''car'', ''driver'', ''ferrari'', ''qualifying'', ''prix'', ''session'', ''fastest'', ''ahead'', ''caution'', ''hamilton'']'
$ws.Range("C27").Value2 = 0.899589057229355
$ws.Range("D27").Value2 = 0.898335132169527
$ws.Range("E27").Value2 = 0.8977043810937559
$ws.Range("F27").Value2 = 0.8974674594264279
$ws.Range("G27").Value2 = 0.8950169921998243

$ws.Range("B28").Value2 = '[''yard'', ''yards'', ''touchdown'', ''tech'', ''alabama'', ''bowl'', ''quarter'', ''michigan'', ''offense'', ''football'', ''pass'', ''rushing'', ''quarterback'', ''conference'', ''touchdowns'']'
$ws.Range("C28").Value2 = 0.894528845283853
$ws.Range("D28").Value2 = 0.8942743768305161
$ws.Range("E28").Value2 = 0.8915809848785308
$ws.Range("F28").Value2 = 0.8822598376371611
$ws.Range("G28").Value2 = 0.8818717964523473

$ws.Range("B29").Value2 = '[''wine'', ''chicken'', ''cheese'', ''bacon'', ''recipes'', ''sandwich'', ''fried'', ''dish'', ''dishes'', ''cuisine'', ''ingredients'', ''cooking'', ''beef'', ''food'', ''product'']'
$ws.Range("C29").Value2 = 0.8959564721457992
$ws.Range("D29").Value2 = 0.8764875153097185
$ws.Range("E29").Value2 = 0.8751223183597634
$ws.Range("F29").Value2 = 0.8695782306780427
$ws.Range("G29").Value2 = 0.8405835638130674

$ws.Range("B30").Value2 = '[''persian'', ''army'', ''byzantine'', ''athens'', ''greeks'', ''greece'', ''battle'', ''alexander'', ''greek'', ''siege'', ''muslim'', ''cavalry'', ''byzantines'', ''arab'', ''ottoman'']'
$ws.Range("C30").Value2 = 0.8918801717134893
$ws.Range("D30").Value2 = 0.8818350529814225
$ws.Range("E30").Value2 = 0.8652987059208482
$ws.Range("F30").Value2 = 0.8382668064278899
$ws.Range("G30").Value2 = 0.8361667774863114

$ws.Range("B31").Value2 = '[''bach'', ''text'', ''movements'', ''movement'', ''soprano'', ''gospel'', ''jesus'', ''hebrew'', ''aria'', ''alto'', ''manuscripts'', ''leipzig'', ''tenor'', ''hymn'', ''matthew'']'
$ws.Range("C31").Value2 = 0.8252006431725999
$ws.Range("D31").Value2 = 0.7995705958719048
$ws.Range("E31").Value2 = 0.7897658468010489
$ws.Range("F31").Value2 = 0.694089658257831
$ws.Range("G31").Value2 = 0.6796103510471005

$ws.Range("B32").Value2 = '[''coaster'', ''ride'', ''roller'', ''riders'', ''train'', ''coasters'', ''flags'', ''park'', ''lift'', ''cedar'', ''drop'', ''brake'', ''trains'', ''steel'', ''hill'']'
$ws.Range("C32").Value2 = 0.8957994125240196
$ws.Range("D32").Value2 = 0.8944680211075516
$ws.Range("E32").Value2 = 0.882740977230002
$ws.Range("F32").Value2 = 0.8788186276840939
$ws.Range("G32").Value2 = 0.8685395620705895

$ws.Range("B33").Value2 = '[''formula'', ''function'', ''matrix'', ''linear'', ''functions'', ''space'', ''constant'', ''defined'', ''language'', ''filter'', ''frequency'', ''mass'', ''derivative'', ''kilogram'', ''units'']'
$ws.Range("C33").Value2 = 0.8809363094291861
$ws.Range("D33").Value2 = 0.8727212074205268
$ws.Range("E33").Value2 = 0.8657413410672804
$ws.Range("F33").Value2 = 0.8410292524168366
$ws.Range("G33").Value2 = 0.8278142687388411

$ws.Range("B34").Value2 = '[''amendment'', ''constitution'', ''singapore'', ''law'', ''court'', ''courts'', ''judicial'', ''parliament'', ''article'', ''constitutional'', ''justice'', ''shall'', ''clause'', ''act'', ''rights'']'
$ws.Range("C34").Value2 = 0.896567062441377
$ws.Range("D34").Value2 = 0.856853908291624
$ws.Range("E34").Value2 = 0.8499817690647762
$ws.Range("F34").Value2 = 0.847849802827159
$ws.Range("G34").Value2 = 0.8393930103432241

$ws.Range("B35").Value2 = '[''management'', ''twitter'', ''investment'', ''bank'', ''billion'', ''equity'', ''watson'', ''design'', ''banking'', ''company'', ''users'', ''firm'', ''business'', ''merger'', ''assets'']'
$ws.Range("C35").Value2 = 0.7552939216034684
$ws.Range("D35").Value2 = 0.7267929555276822
$ws.Range("E35").Value2 = 0.6967749222492603
$ws.Range("F35").Value2 = 0.6373313443470506
$ws.Range("G35").Value2 = 0.5127592840712791

$ws.Range("B36").Value2 = '[''mosque'', ''temple'', ''congregation'', ''wall'', ''palace'', ''hall'', ''jewish'', ''jews'', ''shrine'', ''tomb'', ''temples'', ''cave'', ''jerusalem'', ''christians'', ''dome'']'
$ws.Range("C36").Value2 = 0.8800210270419596
$ws.Range("D36").Value2 = 0.8603731340226697
$ws.Range("E36").Value2 = 0.8394945352111849
$ws.Range("F36").Value2 = 0.8389417350647215
$ws.Range("G36").Value2 = 0.8371138017460441

$ws.Range("B37").Value2 = '[''scotland'', ''scottish'', ''island'', ''islands'', ''norse'', ''century'', ''edinburgh'', ''isles'', ''houses'', ''architecture'', ''mainland'', ''scots'', ''glasgow'', ''highlands'', ''churches'']'
$ws.Range("C37").Value2 = 0.8927239821043242
$ws.Range("D37").Value2 = 0.8429310496340747
$ws.Range("E37").Value2 = 0.8420592345831623
$ws.Range("F37").Value2 = 0.8257043321464712
$ws.Range("G37").Value2 = 0.8000422618234712

$ws.Range("B38").Value2 = '[''croatian'', ''yugoslav'', ''partisans'', ''serbian'', ''croatia'', ''yugoslavia'', ''partisan'', ''serbs'', ''bosnia'', ''belgrade'', ''serbia'', ''germans'', ''civilians'', ''division'', ''forces'']'
$ws.Range("C38").Value2 = 0.866186533677501
$ws.Range("D38").Value2 = 0.8639887463519905
$ws.Range("E38").Value2 = 0.7779175456682912
$ws.Range("F38").Value2 = 0.7770886357257731
$ws.Range("G38").Value2 = 0.7506051997966334

$ws.Range("B39").Value2 = '[''plants'', ''plant'', ''botanical'', ''tree'', ''camouflage'', ''animals'', ''organisms'', ''animal'', ''gardens'', ''species'', ''garden'', ''tea'', ''trees'', ''ecology'', ''predators'']'
$ws.Range("C39").Value2 = 0.8451232096314275
$ws.Range("D39").Value2 = 0.815992496019847
$ws.Range("E39").Value2 = 0.8104518133702522
$ws.Range("F39").Value2 = 0.7960936971394411
$ws.Range("G39").Value2 = 0.7602634450909223

$ws.Range("B40").Value2 = '[''phillies'', ''inning'', ''yankees'', ''dodgers'', ''teams'', ''breaker'', ''giants'', ''mlb'', ''tie'', ''yankee'', ''postseason'', ''game'', ''run'', ''pitcher'', ''baseball'']'
$ws.Range("C40").Value2 = 0.8994062182118583
$ws.Range("D40").Value2 = 0.8977955928166224
$ws.Range("E40").Value2 = 0.8945722923500797
$ws.Range("F40").Value2 = 0.8936233587903845
$ws.Range("G40").Value2 = 0.8931667805655698

$ws.Range("B41").Value2 = '[''painting'', ''paintings'', ''art'', ''artist'', ''serbian'', ''works'', ''marie'', ''exhibition'', ''copenhagen'', ''photography'', ''photographers'', ''canvas'', ''painted'', ''symphony'', ''portrait'']'
$ws.Range("C41").Value2 = 0.8986584859247188
$ws.Range("D41").Value2 = 0.895602751561582
$ws.Range("E41").Value2 = 0.8872227294448619
$ws.Range("F41").Value2 = 0.8811933738725986
$ws.Range("G41").Value2 = 0.8690559314362738

$ws.Range("B42").Value2 = '[''spacecraft'', ''apollo'', ''nasa'', ''orbit'', ''mission'', ''lunar'', ''launch'', ''saturn'', ''docking'', ''flight'', ''manned'', ''module'', ''space'', ''landing'', ''earth'']'
$ws.Range("C42").Value2 = 0.869363424396072
$ws.Range("D42").Value2 = 0.8639757808198418
$ws.Range("E42").Value2 = 0.8611228103650478
$ws.Range("F42").Value2 = 0.8599091447790976
$ws.Range("G42").Value2 = 0.8569662011247938

$ws.Range("B43").Value2 = '[''contest'', ''broadcaster'', ''semi'', ''countries'', ''jury'', ''final'', ''participating'', ''greece'', ''voting'', ''host'', ''idol'', ''entry'', ''philippine'', ''song'', ''country'']'
$ws.Range("C43").Value2 = 0.8977674172528227
$ws.Range("D43").Value2 = 0.897027332936389
$ws.Range("E43").Value2 = 0.8914754394199775
$ws.Range("F43").Value2 = 0.8882109899629327
$ws.Range("G43").Value2 = 0.8860614465100162

$ws.Range("B44").Value2 = '[''lighthouse'', ''light'', ''keeper'', ''tower'', ''lens'', ''keepers'', ''concrete'', ''connecticut'', ''lamp'', ''cottages'', ''installed'', ''constructed'', ''lamps'', ''island'', ''iron'']'
$ws.Range("C44").Value2 = 0.8892793712174863
$ws.Range("D44").Value2 = 0.8758949069477608
$ws.Range("E44").Value2 = 0.8724615980009904
$ws.Range("F44").Value2 = 0.865858147221465
$ws.Range("G44").Value2 = 0.8555457530311984

$ws.Range("B45").Value2 = '[''ben'', ''survivors'', ''shannon'', ''island'', ''freighter'', ''charlotte'', ''bernard'', ''michael'', ''plane'', ''daniel'', ''tom'', ''frank'', ''dave'', ''kate'', ''oceanic'']'
$ws.Range("C45").Value2 = 0.8843045125965368
$ws.Range("D45").Value2 = 0.8697161586874406
$ws.Range("E45").Value2 = 0.8588155624058443
$ws.Range("F45").Value2 = 0.842997492066711
$ws.Range("G45").Value2 = 0.8361383388954475

$ws.Range("B46").Value2 = '[''motorway'', ''croatia'', ''adriatic'', ''croatian'', ''toll'', ''traffic'', ''interchanges'', ''kilometre'', ''route'', ''interchange'', ''kilometres'', ''section'', ''areas'', ''rest'', ''yugoslavia'']'
$ws.Range("C46").Value2 = 0.8974268013124268
$ws.Range("D46").Value2 = 0.8937423833377387
$ws.Range("E46").Value2 = 0.8931426754495102
$ws.Range("F46").Value2 = 0.8849775354748642
$ws.Range("G46").Value2 = 0.8681855529047939

$ws.Range("B47").Value2 = '[''rockets'', ''nba'', ''arena'', ''playoffs'', ''houston'', ''team'', ''franchise'', ''ownership'', ''finals'', ''relocation'', ''toronto'', ''draft'', ''games'', ''pick'', ''season'']'
$ws.Range("C47").Value2 = 0.8994844352280338
$ws.Range("D47").Value2 = 0.8965042202987602
$ws.Range("E47").Value2 = 0.8829335302376713
$ws.Range("F47").Value2 = 0.8826075439356167
$ws.Range("G47").Value2 = 0.8796285845221472

$ws.Range("B48").Value2 = '[''clark'', ''superman'', ''oliver'', ''finale'', ''season'', ''comic'', ''character'', ''relationship'', ''whitney'', ''metropolis'', ''believes'', ''discovers'', ''secret'', ''martha'', ''series'']'
$ws.Range("C48").Value2 = 0.8996513236405467
$ws.Range("D48").Value2 = 0.8995021322793546
$ws.Range("E48").Value2 = 0.8878878148021422
$ws.Range("F48").Value2 = 0.850146904005964
$ws.Range("G48").Value2 = 0.8500225861510472

$ws.Range("B49").Value2 = '[''grammy'', ''neo'', ''hawaiian'', ''soul'', ''nominees'', ''category'', ''awards'', ''rap'', ''award'', ''presented'', ''categories'', ''artists'', ''recipients'', ''academy'', ''disco'']'
$ws.Range("C49").Value2 = 0.899989502367133
$ws.Range("D49").Value2 = 0.899732008444733
$ws.Range("E49").Value2 = 0.8992924662338843
$ws.Range("F49").Value2 = 0.8990640910243387
$ws.Range("G49").Value2 = 0.8987665966169236

$ws.Range("B50").Value2 = '[''euro'', ''coins'', ''currency'', ''note'', ''notes'', ''denominations'', ''dollar'', ''value'', ''stripe'', ''tender'', ''silver'', ''thread'', ''ink'', ''issued'', ''signature'']'
$ws.Range("C50").Value2 = 0.8772465903637001
$ws.Range("D50").Value2 = 0.8701030949011572
$ws.Range("E50").Value2 = 0.864386230131906
$ws.Range("F50").Value2 = 0.8626395032219362
$ws.Range("G50").Value2 = 0.8521718393378128

# Delete row 51 (was id_tpc 49), shrinking the sheet to A1:G50
$ws.Rows(51).Delete()